$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.761.64'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.651.20'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.05'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0631'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.37'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.880.09'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.667.57'
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.22'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.536'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.82'
$ws.Range('E16').Value = '  +5.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.790.84'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0758'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.73'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.40'
$ws.Range('E22').Value = '  +2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.60'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  +13.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.71'
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.12'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.01'
$ws.Range('E29').Value = '  +3.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0525'
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('E32').Value = '  +4.56%  '
$ws.Range('E33').Value = '  +4.50%  '
$ws.Range('E34').Value = '  +4.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.298.67'
$ws.Range('E35').Value = '  +9.32%  '
$ws.Range('E36').Value = '  +5.56%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.835'
$ws.Range('E38').Value = '  +3.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.529'
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.814'
$ws.Range('E41').Value = '  +2.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('E42').Value = '  -3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.45'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.791.76'
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '94.10'
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.16'
$ws.Range('E46').Value = '  +11.53%  '
$ws.Range('E47').Value = '  +5.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.85'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0982'
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.408'
$ws.Range('E51').Value = '  -0.76%  '
